$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that previously held the "Resolving-Mac" sending-cluster block
# (rows 14-17), since that block no longer exists in the updated data.
$ws.Rows("14:17").Delete()

# Updated data for rows 2-13 (recomputed TPM-based NATMI values).
$data = @{
    2 = @("ECs", "Sema3d", "Nrp1", "ECs", 1, 0.3333333333333333, 0.01662366666666667, 0.049871, 0.009202902869893639, 0.00920290286989364, 3, 1, 127.3992563333333, 382.197769, 0.4838549810199306, 0.4838549810199307, 2.117842770866555, 19.060584937799, 0.004452870393440652, 0.004452870393440654)
    3 = @("ECs", "Sema3d", "Nrp1", "FAPs", 1, 0.3333333333333333, 0.01662366666666667, 0.049871, 0.009202902869893639, 0.00920290286989364, 3, 1, 59.36586533333332, 178.097596, 0.2254681108101269, 0.2254681108101269, 0.9868783566795554, 8.881905210115999, 0.002074961124044014, 0.002074961124044014)
    4 = @("ECs", "Sema3d", "Nrp1", "MuSCs", 1, 0.3333333333333333, 0.01662366666666667, 0.049871, 0.009202902869893639, 0.00920290286989364, 3, 1, 16.63275166666667, 49.898255, 0.06317022542837675, 0.06317022542837675, 0.2764973194561111, 2.488475875105, 0.0005813494488866365, 0.0005813494488866366)
    5 = @("ECs", "Sema3d", "Nrp1", "Resolving-Mac", 1, 0.3333333333333333, 0.01662366666666667, 0.049871, 0.009202902869893639, 0.00920290286989364, 3, 1, 59.90262233333334, 179.707867, 0.2275066827415657, 0.2275066827415658, 0.9958012261285556, 8.962211035157, 0.002093721903522337, 0.002093721903522338)
    6 = @("FAPs", "Sema3d", "Nrp1", "ECs", 3, 1, 1.381952, 4.145856, 0.7650520358638441, 0.7650520358638442, 3, 1, 127.3992563333333, 382.197769, 0.4838549810199306, 0.4838549810199307, 176.0596570883627, 1584.536913795264, 0.3701742382921596, 0.3701742382921597)
    7 = @("FAPs", "Sema3d", "Nrp1", "FAPs", 3, 1, 1.381952, 4.145856, 0.7650520358638441, 0.7650520358638442, 3, 1, 59.36586533333332, 178.097596, 0.2254681108101269, 0.2254681108101269, 82.04077632913067, 738.3669869621759, 0.1724948371976623, 0.1724948371976624)
    8 = @("FAPs", "Sema3d", "Nrp1", "MuSCs", 3, 1, 1.381952, 4.145856, 0.7650520358638441, 0.7650520358638442, 3, 1, 16.63275166666667, 49.898255, 0.06317022542837675, 0.06317022542837675, 22.98566443125334, 206.87097988128, 0.0483285095699576, 0.04832850956995761)
    9 = @("FAPs", "Sema3d", "Nrp1", "Resolving-Mac", 3, 1, 1.381952, 4.145856, 0.7650520358638441, 0.7650520358638442, 3, 1, 59.90262233333334, 179.707867, 0.2275066827415657, 0.2275066827415658, 82.78254873879467, 745.0429386491521, 0.1740544508040645, 0.1740544508040646)
    10 = @("MuSCs", "Sema3d", "Nrp1", "ECs", 3, 1, 0.4077746666666667, 1.223324, 0.2257450612662623, 0.2257450612662623, 3, 1, 127.3992563333333, 382.197769, 0.4838549810199306, 0.4838549810199307, 51.95018928490622, 467.551703564156, 0.1092278723343304, 0.1092278723343305)
    11 = @("MuSCs", "Sema3d", "Nrp1", "FAPs", 3, 1, 0.4077746666666667, 1.223324, 0.2257450612662623, 0.2257450612662623, 3, 1, 59.36586533333332, 178.097596, 0.2254681108101269, 0.2254681108101269, 24.20789594767822, 217.871063529104, 0.05089831248842051, 0.05089831248842053)
    12 = @("MuSCs", "Sema3d", "Nrp1", "MuSCs", 3, 1, 0.4077746666666667, 1.223324, 0.2257450612662623, 0.2257450612662623, 3, 1, 16.63275166666667, 49.898255, 0.06317022542837675, 0.06317022542837675, 6.782414766624445, 61.04173289962, 0.01426036640953251, 0.01426036640953251)
    13 = @("MuSCs", "Sema3d", "Nrp1", "Resolving-Mac", 3, 1, 0.4077746666666667, 1.223324, 0.2257450612662623, 0.2257450612662623, 3, 1, 59.90262233333334, 179.707867, 0.2275066827415657, 0.2275066827415658, 24.42677185443423, 219.840946689908, 0.05135851003397886, 0.05135851003397888)
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}